$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -4
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = -3
$ws.Range("F7").Value = 4
$ws.Range("F8").Value = -1
$ws.Range("F9").Value = 4
$ws.Range("F10").Value = -4
$ws.Range("F11").Value = -6
